# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (used by the notes master)
#   ppt/theme/theme2.xml  -> "Integral" / "Red Violet" (used by the slide master,
#                             i.e. the theme that actually paints the slides)
#
# The authored edit swaps the two themes' contents: the slide master's theme
# becomes the plain "Office Theme" colour scheme (and the notes master's
# theme becomes "Integral"). The only thing that differs between the two
# theme parts is the <a:clrScheme> (font/format schemes are identical), so
# we reproduce the visible half of that swap - the slide master's colour
# scheme - via the ThemeColorScheme COM surface, which is the one that
# actually round-trips cleanly back into ppt/theme/theme2.xml.

function ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# "Office Theme" colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $tcs.Colors($i).RGB = ComRGB($officeColors[$i - 1])
}
